$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132846713066101
$ws.Range("B1").Value = 2.252492189407349
$ws.Range("C1").Value = 10.94107532501221
$ws.Range("D1").Value = 2.206290006637573
$ws.Range("E1").Value = 1.28030788898468
